# Trade #58 closed at 2026-02-17 15:43:42 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status aggregate rows and appends the new
# trade row (#58) to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.38
$wsSummary.Range("B4").Value = 0.38
$wsSummary.Range("B5").Value = 0.13
$wsSummary.Range("B6").Value = 58
$wsSummary.Range("B7").Value = 18
$wsSummary.Range("B9").Value = 31.03

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.38
$wsStatus.Range("D4").Value = 58
$wsStatus.Range("E4").Value = 0.38
$wsStatus.Range("F4").Value = 0.38
$wsStatus.Range("G4").Value = 31.03

# ---------------------------------------------------------------------
# 3. Append the new trade (#58) to a trades-log worksheet
# ---------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Cells.Item(59, 1).Value = 58
    # Force text storage for the date string so Excel does not auto-coerce
    # "2026-02-17" into a date serial number (matches the plain text cells
    # used for every other date/time/text field in this sheet), then drop
    # back to the default "Normal" style so no stray formatting is left
    # behind on the cell.
    $ws.Cells.Item(59, 2).NumberFormat = "@"
    $ws.Cells.Item(59, 2).Value = "2026-02-17"
    $ws.Cells.Item(59, 2).Style = "Normal"
    $ws.Cells.Item(59, 3).Value = "15:43:35"
    $ws.Cells.Item(59, 4).Value = "MarketMaking"
    $ws.Cells.Item(59, 5).Value = "UP"
    $ws.Cells.Item(59, 6).Value = 0.15
    $ws.Cells.Item(59, 7).Value = 0.182692
    $ws.Cells.Item(59, 8).Value = "CLOSED"
    $ws.Cells.Item(59, 9).Value = 21.7949
    $ws.Cells.Item(59, 10).Value = 0.03
    $ws.Cells.Item(59, 11).Value = 100.38
    $ws.Cells.Item(59, 12).Value = 0
    $ws.Cells.Item(59, 13).Value = 0
    $ws.Cells.Item(59, 14).Value = 0.6
    $ws.Cells.Item(59, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(59, 16).Value = "early_exit"
    $ws.Cells.Item(59, 17).Value = 0.15
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
